$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 1: numeric column indices 0-11 (keeps existing bold/border style s=1)
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10
$ws.Cells.Item(1, 12).Value = 11

# Row 2: old header text labels (default style, no border/bold); I2, K2, L2 cleared
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "Lg."
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "Threading"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "HeadDia."
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "Head Ht."
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "DriveSize"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "Tensile Strength"
$ws.Cells.Item(2, 6).Style = "Normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "Specifications Met"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "Pkg.Qty."
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 9).Value = $null
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "Pkg."
$ws.Cells.Item(2, 10).Style = "Normal"
$ws.Cells.Item(2, 11).Value = $null
$ws.Cells.Item(2, 12).Value = $null

# Rows 3-39: old data rows 2-38 shifted down by one (forced as text to match source formatting)
# row 3
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "3/16`""
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "Fully Threaded"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "0.167`""
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "0.062`""
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "T8"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "Not Rated"
$ws.Cells.Item(3, 6).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = "25"
$ws.Cells.Item(3, 8).Style = "Normal"
$ws.Cells.Item(3, 9).NumberFormat = "@"
$ws.Cells.Item(3, 9).Value = "97690A142"
$ws.Cells.Item(3, 9).Style = "Normal"
$ws.Cells.Item(3, 10).NumberFormat = "@"
$ws.Cells.Item(3, 10).Value = "`$6.71"
$ws.Cells.Item(3, 10).Style = "Normal"
$ws.Cells.Item(3, 11).NumberFormat = "@"
$ws.Cells.Item(3, 11).Value = "2-56"
$ws.Cells.Item(3, 11).Style = "Normal"
$ws.Cells.Item(3, 12).NumberFormat = "@"
$ws.Cells.Item(3, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(3, 12).Style = "Normal"
# row 4
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "1/4`""
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "Fully Threaded"
$ws.Cells.Item(4, 2).Style = "Normal"
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "0.167`""
$ws.Cells.Item(4, 3).Style = "Normal"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.062`""
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "T8"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "Not Rated"
$ws.Cells.Item(4, 6).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = "25"
$ws.Cells.Item(4, 8).Style = "Normal"
$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = "97690A143"
$ws.Cells.Item(4, 9).Style = "Normal"
$ws.Cells.Item(4, 10).NumberFormat = "@"
$ws.Cells.Item(4, 10).Value = "7.15"
$ws.Cells.Item(4, 10).Style = "Normal"
$ws.Cells.Item(4, 11).NumberFormat = "@"
$ws.Cells.Item(4, 11).Value = "2-56"
$ws.Cells.Item(4, 11).Style = "Normal"
$ws.Cells.Item(4, 12).NumberFormat = "@"
$ws.Cells.Item(4, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(4, 12).Style = "Normal"
# row 5
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "5/16`""
$ws.Cells.Item(5, 1).Style = "Normal"
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "Fully Threaded"
$ws.Cells.Item(5, 2).Style = "Normal"
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "0.167`""
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.062`""
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "T8"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "Not Rated"
$ws.Cells.Item(5, 6).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value = "25"
$ws.Cells.Item(5, 8).Style = "Normal"
$ws.Cells.Item(5, 9).NumberFormat = "@"
$ws.Cells.Item(5, 9).Value = "97690A144"
$ws.Cells.Item(5, 9).Style = "Normal"
$ws.Cells.Item(5, 10).NumberFormat = "@"
$ws.Cells.Item(5, 10).Value = "6.82"
$ws.Cells.Item(5, 10).Style = "Normal"
$ws.Cells.Item(5, 11).NumberFormat = "@"
$ws.Cells.Item(5, 11).Value = "2-56"
$ws.Cells.Item(5, 11).Style = "Normal"
$ws.Cells.Item(5, 12).NumberFormat = "@"
$ws.Cells.Item(5, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(5, 12).Style = "Normal"
# row 6
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "3/8`""
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "Fully Threaded"
$ws.Cells.Item(6, 2).Style = "Normal"
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "0.167`""
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.062`""
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "T8"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "Not Rated"
$ws.Cells.Item(6, 6).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(6, 8).NumberFormat = "@"
$ws.Cells.Item(6, 8).Value = "25"
$ws.Cells.Item(6, 8).Style = "Normal"
$ws.Cells.Item(6, 9).NumberFormat = "@"
$ws.Cells.Item(6, 9).Value = "97690A145"
$ws.Cells.Item(6, 9).Style = "Normal"
$ws.Cells.Item(6, 10).NumberFormat = "@"
$ws.Cells.Item(6, 10).Value = "6.86"
$ws.Cells.Item(6, 10).Style = "Normal"
$ws.Cells.Item(6, 11).NumberFormat = "@"
$ws.Cells.Item(6, 11).Value = "2-56"
$ws.Cells.Item(6, 11).Style = "Normal"
$ws.Cells.Item(6, 12).NumberFormat = "@"
$ws.Cells.Item(6, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(6, 12).Style = "Normal"
# row 7
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "1/2`""
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "Fully Threaded"
$ws.Cells.Item(7, 2).Style = "Normal"
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "0.167`""
$ws.Cells.Item(7, 3).Style = "Normal"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.062`""
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "T8"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = "Not Rated"
$ws.Cells.Item(7, 6).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(7, 8).NumberFormat = "@"
$ws.Cells.Item(7, 8).Value = "50"
$ws.Cells.Item(7, 8).Style = "Normal"
$ws.Cells.Item(7, 9).NumberFormat = "@"
$ws.Cells.Item(7, 9).Value = "97690A146"
$ws.Cells.Item(7, 9).Style = "Normal"
$ws.Cells.Item(7, 10).NumberFormat = "@"
$ws.Cells.Item(7, 10).Value = "9.69"
$ws.Cells.Item(7, 10).Style = "Normal"
$ws.Cells.Item(7, 11).NumberFormat = "@"
$ws.Cells.Item(7, 11).Value = "2-56"
$ws.Cells.Item(7, 11).Style = "Normal"
$ws.Cells.Item(7, 12).NumberFormat = "@"
$ws.Cells.Item(7, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(7, 12).Style = "Normal"
# row 8
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "3/16`""
$ws.Cells.Item(8, 1).Style = "Normal"
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "Fully Threaded"
$ws.Cells.Item(8, 2).Style = "Normal"
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "0.219`""
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.08`""
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "T10"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = "Not Rated"
$ws.Cells.Item(8, 6).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = "25"
$ws.Cells.Item(8, 8).Style = "Normal"
$ws.Cells.Item(8, 9).NumberFormat = "@"
$ws.Cells.Item(8, 9).Value = "97690A147"
$ws.Cells.Item(8, 9).Style = "Normal"
$ws.Cells.Item(8, 10).NumberFormat = "@"
$ws.Cells.Item(8, 10).Value = "5.63"
$ws.Cells.Item(8, 10).Style = "Normal"
$ws.Cells.Item(8, 11).NumberFormat = "@"
$ws.Cells.Item(8, 11).Value = "4-40"
$ws.Cells.Item(8, 11).Style = "Normal"
$ws.Cells.Item(8, 12).NumberFormat = "@"
$ws.Cells.Item(8, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(8, 12).Style = "Normal"
# row 9
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "1/4`""
$ws.Cells.Item(9, 1).Style = "Normal"
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "Fully Threaded"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "0.219`""
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.08`""
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "T10"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = "Not Rated"
$ws.Cells.Item(9, 6).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = "50"
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = "97690A148"
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 10).NumberFormat = "@"
$ws.Cells.Item(9, 10).Value = "9.27"
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(9, 11).NumberFormat = "@"
$ws.Cells.Item(9, 11).Value = "4-40"
$ws.Cells.Item(9, 11).Style = "Normal"
$ws.Cells.Item(9, 12).NumberFormat = "@"
$ws.Cells.Item(9, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(9, 12).Style = "Normal"
# row 10
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "5/16`""
$ws.Cells.Item(10, 1).Style = "Normal"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "Fully Threaded"
$ws.Cells.Item(10, 2).Style = "Normal"
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "0.219`""
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08`""
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "T10"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = "Not Rated"
$ws.Cells.Item(10, 6).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = "50"
$ws.Cells.Item(10, 8).Style = "Normal"
$ws.Cells.Item(10, 9).NumberFormat = "@"
$ws.Cells.Item(10, 9).Value = "97690A149"
$ws.Cells.Item(10, 9).Style = "Normal"
$ws.Cells.Item(10, 10).NumberFormat = "@"
$ws.Cells.Item(10, 10).Value = "9.06"
$ws.Cells.Item(10, 10).Style = "Normal"
$ws.Cells.Item(10, 11).NumberFormat = "@"
$ws.Cells.Item(10, 11).Value = "4-40"
$ws.Cells.Item(10, 11).Style = "Normal"
$ws.Cells.Item(10, 12).NumberFormat = "@"
$ws.Cells.Item(10, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(10, 12).Style = "Normal"
# row 11
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "3/8`""
$ws.Cells.Item(11, 1).Style = "Normal"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "Fully Threaded"
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "0.219`""
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08`""
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "T10"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = "Not Rated"
$ws.Cells.Item(11, 6).Style = "Normal"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = "50"
$ws.Cells.Item(11, 8).Style = "Normal"
$ws.Cells.Item(11, 9).NumberFormat = "@"
$ws.Cells.Item(11, 9).Value = "97690A151"
$ws.Cells.Item(11, 9).Style = "Normal"
$ws.Cells.Item(11, 10).NumberFormat = "@"
$ws.Cells.Item(11, 10).Value = "9.17"
$ws.Cells.Item(11, 10).Style = "Normal"
$ws.Cells.Item(11, 11).NumberFormat = "@"
$ws.Cells.Item(11, 11).Value = "4-40"
$ws.Cells.Item(11, 11).Style = "Normal"
$ws.Cells.Item(11, 12).NumberFormat = "@"
$ws.Cells.Item(11, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(11, 12).Style = "Normal"
# row 12
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "7/16`""
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = "Fully Threaded"
$ws.Cells.Item(12, 2).Style = "Normal"
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "0.219`""
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08`""
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "T10"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = "Not Rated"
$ws.Cells.Item(12, 6).Style = "Normal"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = "50"
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(12, 9).NumberFormat = "@"
$ws.Cells.Item(12, 9).Value = "97690A152"
$ws.Cells.Item(12, 9).Style = "Normal"
$ws.Cells.Item(12, 10).NumberFormat = "@"
$ws.Cells.Item(12, 10).Value = "9.47"
$ws.Cells.Item(12, 10).Style = "Normal"
$ws.Cells.Item(12, 11).NumberFormat = "@"
$ws.Cells.Item(12, 11).Value = "4-40"
$ws.Cells.Item(12, 11).Style = "Normal"
$ws.Cells.Item(12, 12).NumberFormat = "@"
$ws.Cells.Item(12, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(12, 12).Style = "Normal"
# row 13
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "1/2`""
$ws.Cells.Item(13, 1).Style = "Normal"
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "Fully Threaded"
$ws.Cells.Item(13, 2).Style = "Normal"
$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "0.219`""
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.08`""
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "T10"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 6).NumberFormat = "@"
$ws.Cells.Item(13, 6).Value = "Not Rated"
$ws.Cells.Item(13, 6).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(13, 8).NumberFormat = "@"
$ws.Cells.Item(13, 8).Value = "50"
$ws.Cells.Item(13, 8).Style = "Normal"
$ws.Cells.Item(13, 9).NumberFormat = "@"
$ws.Cells.Item(13, 9).Value = "97690A153"
$ws.Cells.Item(13, 9).Style = "Normal"
$ws.Cells.Item(13, 10).NumberFormat = "@"
$ws.Cells.Item(13, 10).Value = "10.19"
$ws.Cells.Item(13, 10).Style = "Normal"
$ws.Cells.Item(13, 11).NumberFormat = "@"
$ws.Cells.Item(13, 11).Value = "4-40"
$ws.Cells.Item(13, 11).Style = "Normal"
$ws.Cells.Item(13, 12).NumberFormat = "@"
$ws.Cells.Item(13, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(13, 12).Style = "Normal"
# row 14
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "3/4`""
$ws.Cells.Item(14, 1).Style = "Normal"
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "Fully Threaded"
$ws.Cells.Item(14, 2).Style = "Normal"
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "0.219`""
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.08`""
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "T10"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = "Not Rated"
$ws.Cells.Item(14, 6).Style = "Normal"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(14, 8).NumberFormat = "@"
$ws.Cells.Item(14, 8).Value = "25"
$ws.Cells.Item(14, 8).Style = "Normal"
$ws.Cells.Item(14, 9).NumberFormat = "@"
$ws.Cells.Item(14, 9).Value = "97690A154"
$ws.Cells.Item(14, 9).Style = "Normal"
$ws.Cells.Item(14, 10).NumberFormat = "@"
$ws.Cells.Item(14, 10).Value = "4.86"
$ws.Cells.Item(14, 10).Style = "Normal"
$ws.Cells.Item(14, 11).NumberFormat = "@"
$ws.Cells.Item(14, 11).Value = "4-40"
$ws.Cells.Item(14, 11).Style = "Normal"
$ws.Cells.Item(14, 12).NumberFormat = "@"
$ws.Cells.Item(14, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(14, 12).Style = "Normal"
# row 15
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "3/16`""
$ws.Cells.Item(15, 1).Style = "Normal"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "Fully Threaded"
$ws.Cells.Item(15, 2).Style = "Normal"
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "0.27`""
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.097`""
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "T15"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 6).NumberFormat = "@"
$ws.Cells.Item(15, 6).Value = "Not Rated"
$ws.Cells.Item(15, 6).Style = "Normal"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "50"
$ws.Cells.Item(15, 8).Style = "Normal"
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = "97690A155"
$ws.Cells.Item(15, 9).Style = "Normal"
$ws.Cells.Item(15, 10).NumberFormat = "@"
$ws.Cells.Item(15, 10).Value = "9.68"
$ws.Cells.Item(15, 10).Style = "Normal"
$ws.Cells.Item(15, 11).NumberFormat = "@"
$ws.Cells.Item(15, 11).Value = "6-32"
$ws.Cells.Item(15, 11).Style = "Normal"
$ws.Cells.Item(15, 12).NumberFormat = "@"
$ws.Cells.Item(15, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(15, 12).Style = "Normal"
# row 16
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "1/4`""
$ws.Cells.Item(16, 1).Style = "Normal"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "Fully Threaded"
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = "0.27`""
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.097`""
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "T15"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 6).NumberFormat = "@"
$ws.Cells.Item(16, 6).Value = "Not Rated"
$ws.Cells.Item(16, 6).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value = "50"
$ws.Cells.Item(16, 8).Style = "Normal"
$ws.Cells.Item(16, 9).NumberFormat = "@"
$ws.Cells.Item(16, 9).Value = "97690A156"
$ws.Cells.Item(16, 9).Style = "Normal"
$ws.Cells.Item(16, 10).NumberFormat = "@"
$ws.Cells.Item(16, 10).Value = "9.79"
$ws.Cells.Item(16, 10).Style = "Normal"
$ws.Cells.Item(16, 11).NumberFormat = "@"
$ws.Cells.Item(16, 11).Value = "6-32"
$ws.Cells.Item(16, 11).Style = "Normal"
$ws.Cells.Item(16, 12).NumberFormat = "@"
$ws.Cells.Item(16, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(16, 12).Style = "Normal"
# row 17
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "5/16`""
$ws.Cells.Item(17, 1).Style = "Normal"
$ws.Cells.Item(17, 2).NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = "Fully Threaded"
$ws.Cells.Item(17, 2).Style = "Normal"
$ws.Cells.Item(17, 3).NumberFormat = "@"
$ws.Cells.Item(17, 3).Value = "0.27`""
$ws.Cells.Item(17, 3).Style = "Normal"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.097`""
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "T15"
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 6).NumberFormat = "@"
$ws.Cells.Item(17, 6).Value = "Not Rated"
$ws.Cells.Item(17, 6).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(17, 8).NumberFormat = "@"
$ws.Cells.Item(17, 8).Value = "50"
$ws.Cells.Item(17, 8).Style = "Normal"
$ws.Cells.Item(17, 9).NumberFormat = "@"
$ws.Cells.Item(17, 9).Value = "97690A157"
$ws.Cells.Item(17, 9).Style = "Normal"
$ws.Cells.Item(17, 10).NumberFormat = "@"
$ws.Cells.Item(17, 10).Value = "9.79"
$ws.Cells.Item(17, 10).Style = "Normal"
$ws.Cells.Item(17, 11).NumberFormat = "@"
$ws.Cells.Item(17, 11).Value = "6-32"
$ws.Cells.Item(17, 11).Style = "Normal"
$ws.Cells.Item(17, 12).NumberFormat = "@"
$ws.Cells.Item(17, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(17, 12).Style = "Normal"
# row 18
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "3/8`""
$ws.Cells.Item(18, 1).Style = "Normal"
$ws.Cells.Item(18, 2).NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = "Fully Threaded"
$ws.Cells.Item(18, 2).Style = "Normal"
$ws.Cells.Item(18, 3).NumberFormat = "@"
$ws.Cells.Item(18, 3).Value = "0.27`""
$ws.Cells.Item(18, 3).Style = "Normal"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.097`""
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "T15"
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 6).NumberFormat = "@"
$ws.Cells.Item(18, 6).Value = "Not Rated"
$ws.Cells.Item(18, 6).Style = "Normal"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value = "50"
$ws.Cells.Item(18, 8).Style = "Normal"
$ws.Cells.Item(18, 9).NumberFormat = "@"
$ws.Cells.Item(18, 9).Value = "97690A158"
$ws.Cells.Item(18, 9).Style = "Normal"
$ws.Cells.Item(18, 10).NumberFormat = "@"
$ws.Cells.Item(18, 10).Value = "9.98"
$ws.Cells.Item(18, 10).Style = "Normal"
$ws.Cells.Item(18, 11).NumberFormat = "@"
$ws.Cells.Item(18, 11).Value = "6-32"
$ws.Cells.Item(18, 11).Style = "Normal"
$ws.Cells.Item(18, 12).NumberFormat = "@"
$ws.Cells.Item(18, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(18, 12).Style = "Normal"
# row 19
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "7/16`""
$ws.Cells.Item(19, 1).Style = "Normal"
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = "Fully Threaded"
$ws.Cells.Item(19, 2).Style = "Normal"
$ws.Cells.Item(19, 3).NumberFormat = "@"
$ws.Cells.Item(19, 3).Value = "0.27`""
$ws.Cells.Item(19, 3).Style = "Normal"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.097`""
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "T15"
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 6).NumberFormat = "@"
$ws.Cells.Item(19, 6).Value = "Not Rated"
$ws.Cells.Item(19, 6).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(19, 8).NumberFormat = "@"
$ws.Cells.Item(19, 8).Value = "50"
$ws.Cells.Item(19, 8).Style = "Normal"
$ws.Cells.Item(19, 9).NumberFormat = "@"
$ws.Cells.Item(19, 9).Value = "97690A159"
$ws.Cells.Item(19, 9).Style = "Normal"
$ws.Cells.Item(19, 10).NumberFormat = "@"
$ws.Cells.Item(19, 10).Value = "10.15"
$ws.Cells.Item(19, 10).Style = "Normal"
$ws.Cells.Item(19, 11).NumberFormat = "@"
$ws.Cells.Item(19, 11).Value = "6-32"
$ws.Cells.Item(19, 11).Style = "Normal"
$ws.Cells.Item(19, 12).NumberFormat = "@"
$ws.Cells.Item(19, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(19, 12).Style = "Normal"
# row 20
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "1/2`""
$ws.Cells.Item(20, 1).Style = "Normal"
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = "Fully Threaded"
$ws.Cells.Item(20, 2).Style = "Normal"
$ws.Cells.Item(20, 3).NumberFormat = "@"
$ws.Cells.Item(20, 3).Value = "0.27`""
$ws.Cells.Item(20, 3).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.097`""
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "T15"
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 6).NumberFormat = "@"
$ws.Cells.Item(20, 6).Value = "Not Rated"
$ws.Cells.Item(20, 6).Style = "Normal"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(20, 8).NumberFormat = "@"
$ws.Cells.Item(20, 8).Value = "50"
$ws.Cells.Item(20, 8).Style = "Normal"
$ws.Cells.Item(20, 9).NumberFormat = "@"
$ws.Cells.Item(20, 9).Value = "97690A161"
$ws.Cells.Item(20, 9).Style = "Normal"
$ws.Cells.Item(20, 10).NumberFormat = "@"
$ws.Cells.Item(20, 10).Value = "10.45"
$ws.Cells.Item(20, 10).Style = "Normal"
$ws.Cells.Item(20, 11).NumberFormat = "@"
$ws.Cells.Item(20, 11).Value = "6-32"
$ws.Cells.Item(20, 11).Style = "Normal"
$ws.Cells.Item(20, 12).NumberFormat = "@"
$ws.Cells.Item(20, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(20, 12).Style = "Normal"
# row 21
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "3/4`""
$ws.Cells.Item(21, 1).Style = "Normal"
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = "Fully Threaded"
$ws.Cells.Item(21, 2).Style = "Normal"
$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "0.27`""
$ws.Cells.Item(21, 3).Style = "Normal"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.097`""
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "T15"
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 6).NumberFormat = "@"
$ws.Cells.Item(21, 6).Value = "Not Rated"
$ws.Cells.Item(21, 6).Style = "Normal"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(21, 8).NumberFormat = "@"
$ws.Cells.Item(21, 8).Value = "25"
$ws.Cells.Item(21, 8).Style = "Normal"
$ws.Cells.Item(21, 9).NumberFormat = "@"
$ws.Cells.Item(21, 9).Value = "97690A162"
$ws.Cells.Item(21, 9).Style = "Normal"
$ws.Cells.Item(21, 10).NumberFormat = "@"
$ws.Cells.Item(21, 10).Value = "5.55"
$ws.Cells.Item(21, 10).Style = "Normal"
$ws.Cells.Item(21, 11).NumberFormat = "@"
$ws.Cells.Item(21, 11).Value = "6-32"
$ws.Cells.Item(21, 11).Style = "Normal"
$ws.Cells.Item(21, 12).NumberFormat = "@"
$ws.Cells.Item(21, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(21, 12).Style = "Normal"
# row 22
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "1`""
$ws.Cells.Item(22, 1).Style = "Normal"
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = "Fully Threaded"
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 3).NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = "0.27`""
$ws.Cells.Item(22, 3).Style = "Normal"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.097`""
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "T15"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 6).NumberFormat = "@"
$ws.Cells.Item(22, 6).Value = "Not Rated"
$ws.Cells.Item(22, 6).Style = "Normal"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(22, 8).NumberFormat = "@"
$ws.Cells.Item(22, 8).Value = "25"
$ws.Cells.Item(22, 8).Style = "Normal"
$ws.Cells.Item(22, 9).NumberFormat = "@"
$ws.Cells.Item(22, 9).Value = "97690A163"
$ws.Cells.Item(22, 9).Style = "Normal"
$ws.Cells.Item(22, 10).NumberFormat = "@"
$ws.Cells.Item(22, 10).Value = "6.46"
$ws.Cells.Item(22, 10).Style = "Normal"
$ws.Cells.Item(22, 11).NumberFormat = "@"
$ws.Cells.Item(22, 11).Value = "6-32"
$ws.Cells.Item(22, 11).Style = "Normal"
$ws.Cells.Item(22, 12).NumberFormat = "@"
$ws.Cells.Item(22, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(22, 12).Style = "Normal"
# row 23
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "3/16`""
$ws.Cells.Item(23, 1).Style = "Normal"
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = "Fully Threaded"
$ws.Cells.Item(23, 2).Style = "Normal"
$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = "0.322`""
$ws.Cells.Item(23, 3).Style = "Normal"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.115`""
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "T20"
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 6).NumberFormat = "@"
$ws.Cells.Item(23, 6).Value = "Not Rated"
$ws.Cells.Item(23, 6).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(23, 8).NumberFormat = "@"
$ws.Cells.Item(23, 8).Value = "50"
$ws.Cells.Item(23, 8).Style = "Normal"
$ws.Cells.Item(23, 9).NumberFormat = "@"
$ws.Cells.Item(23, 9).Value = "97690A164"
$ws.Cells.Item(23, 9).Style = "Normal"
$ws.Cells.Item(23, 10).NumberFormat = "@"
$ws.Cells.Item(23, 10).Value = "8.70"
$ws.Cells.Item(23, 10).Style = "Normal"
$ws.Cells.Item(23, 11).NumberFormat = "@"
$ws.Cells.Item(23, 11).Value = "8-32"
$ws.Cells.Item(23, 11).Style = "Normal"
$ws.Cells.Item(23, 12).NumberFormat = "@"
$ws.Cells.Item(23, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(23, 12).Style = "Normal"
# row 24
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "1/4`""
$ws.Cells.Item(24, 1).Style = "Normal"
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = "Fully Threaded"
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 3).NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = "0.322`""
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.115`""
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "T20"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 6).NumberFormat = "@"
$ws.Cells.Item(24, 6).Value = "Not Rated"
$ws.Cells.Item(24, 6).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(24, 8).NumberFormat = "@"
$ws.Cells.Item(24, 8).Value = "50"
$ws.Cells.Item(24, 8).Style = "Normal"
$ws.Cells.Item(24, 9).NumberFormat = "@"
$ws.Cells.Item(24, 9).Value = "97690A165"
$ws.Cells.Item(24, 9).Style = "Normal"
$ws.Cells.Item(24, 10).NumberFormat = "@"
$ws.Cells.Item(24, 10).Value = "10.58"
$ws.Cells.Item(24, 10).Style = "Normal"
$ws.Cells.Item(24, 11).NumberFormat = "@"
$ws.Cells.Item(24, 11).Value = "8-32"
$ws.Cells.Item(24, 11).Style = "Normal"
$ws.Cells.Item(24, 12).NumberFormat = "@"
$ws.Cells.Item(24, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(24, 12).Style = "Normal"
# row 25
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "5/16`""
$ws.Cells.Item(25, 1).Style = "Normal"
$ws.Cells.Item(25, 2).NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = "Fully Threaded"
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = "0.322`""
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.115`""
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "T20"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 6).NumberFormat = "@"
$ws.Cells.Item(25, 6).Value = "Not Rated"
$ws.Cells.Item(25, 6).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = "25"
$ws.Cells.Item(25, 8).Style = "Normal"
$ws.Cells.Item(25, 9).NumberFormat = "@"
$ws.Cells.Item(25, 9).Value = "97690A166"
$ws.Cells.Item(25, 9).Style = "Normal"
$ws.Cells.Item(25, 10).NumberFormat = "@"
$ws.Cells.Item(25, 10).Value = "5.88"
$ws.Cells.Item(25, 10).Style = "Normal"
$ws.Cells.Item(25, 11).NumberFormat = "@"
$ws.Cells.Item(25, 11).Value = "8-32"
$ws.Cells.Item(25, 11).Style = "Normal"
$ws.Cells.Item(25, 12).NumberFormat = "@"
$ws.Cells.Item(25, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(25, 12).Style = "Normal"
# row 26
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "3/8`""
$ws.Cells.Item(26, 1).Style = "Normal"
$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "Fully Threaded"
$ws.Cells.Item(26, 2).Style = "Normal"
$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "0.322`""
$ws.Cells.Item(26, 3).Style = "Normal"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.115`""
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "T20"
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 6).NumberFormat = "@"
$ws.Cells.Item(26, 6).Value = "Not Rated"
$ws.Cells.Item(26, 6).Style = "Normal"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "50"
$ws.Cells.Item(26, 8).Style = "Normal"
$ws.Cells.Item(26, 9).NumberFormat = "@"
$ws.Cells.Item(26, 9).Value = "97690A167"
$ws.Cells.Item(26, 9).Style = "Normal"
$ws.Cells.Item(26, 10).NumberFormat = "@"
$ws.Cells.Item(26, 10).Value = "10.68"
$ws.Cells.Item(26, 10).Style = "Normal"
$ws.Cells.Item(26, 11).NumberFormat = "@"
$ws.Cells.Item(26, 11).Value = "8-32"
$ws.Cells.Item(26, 11).Style = "Normal"
$ws.Cells.Item(26, 12).NumberFormat = "@"
$ws.Cells.Item(26, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(26, 12).Style = "Normal"
# row 27
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "7/16`""
$ws.Cells.Item(27, 1).Style = "Normal"
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "Fully Threaded"
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "0.322`""
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.115`""
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "T20"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = "Not Rated"
$ws.Cells.Item(27, 6).Style = "Normal"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value = "50"
$ws.Cells.Item(27, 8).Style = "Normal"
$ws.Cells.Item(27, 9).NumberFormat = "@"
$ws.Cells.Item(27, 9).Value = "97690A168"
$ws.Cells.Item(27, 9).Style = "Normal"
$ws.Cells.Item(27, 10).NumberFormat = "@"
$ws.Cells.Item(27, 10).Value = "11.49"
$ws.Cells.Item(27, 10).Style = "Normal"
$ws.Cells.Item(27, 11).NumberFormat = "@"
$ws.Cells.Item(27, 11).Value = "8-32"
$ws.Cells.Item(27, 11).Style = "Normal"
$ws.Cells.Item(27, 12).NumberFormat = "@"
$ws.Cells.Item(27, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(27, 12).Style = "Normal"
# row 28
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "1/2`""
$ws.Cells.Item(28, 1).Style = "Normal"
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = "Fully Threaded"
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "0.322`""
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.115`""
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "T20"
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(28, 6).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value = "Not Rated"
$ws.Cells.Item(28, 6).Style = "Normal"
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 8).Value = "25"
$ws.Cells.Item(28, 8).Style = "Normal"
$ws.Cells.Item(28, 9).NumberFormat = "@"
$ws.Cells.Item(28, 9).Value = "97690A169"
$ws.Cells.Item(28, 9).Style = "Normal"
$ws.Cells.Item(28, 10).NumberFormat = "@"
$ws.Cells.Item(28, 10).Value = "5.58"
$ws.Cells.Item(28, 10).Style = "Normal"
$ws.Cells.Item(28, 11).NumberFormat = "@"
$ws.Cells.Item(28, 11).Value = "8-32"
$ws.Cells.Item(28, 11).Style = "Normal"
$ws.Cells.Item(28, 12).NumberFormat = "@"
$ws.Cells.Item(28, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(28, 12).Style = "Normal"
# row 29
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "3/4`""
$ws.Cells.Item(29, 1).Style = "Normal"
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = "Fully Threaded"
$ws.Cells.Item(29, 2).Style = "Normal"
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "0.322`""
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.115`""
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "T20"
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(29, 6).NumberFormat = "@"
$ws.Cells.Item(29, 6).Value = "Not Rated"
$ws.Cells.Item(29, 6).Style = "Normal"
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(29, 8).NumberFormat = "@"
$ws.Cells.Item(29, 8).Value = "25"
$ws.Cells.Item(29, 8).Style = "Normal"
$ws.Cells.Item(29, 9).NumberFormat = "@"
$ws.Cells.Item(29, 9).Value = "97690A171"
$ws.Cells.Item(29, 9).Style = "Normal"
$ws.Cells.Item(29, 10).NumberFormat = "@"
$ws.Cells.Item(29, 10).Value = "6.48"
$ws.Cells.Item(29, 10).Style = "Normal"
$ws.Cells.Item(29, 11).NumberFormat = "@"
$ws.Cells.Item(29, 11).Value = "8-32"
$ws.Cells.Item(29, 11).Style = "Normal"
$ws.Cells.Item(29, 12).NumberFormat = "@"
$ws.Cells.Item(29, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(29, 12).Style = "Normal"
# row 30
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "1`""
$ws.Cells.Item(30, 1).Style = "Normal"
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "Fully Threaded"
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "0.322`""
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.115`""
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "T20"
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = "Not Rated"
$ws.Cells.Item(30, 6).Style = "Normal"
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(30, 8).NumberFormat = "@"
$ws.Cells.Item(30, 8).Value = "25"
$ws.Cells.Item(30, 8).Style = "Normal"
$ws.Cells.Item(30, 9).NumberFormat = "@"
$ws.Cells.Item(30, 9).Value = "97690A172"
$ws.Cells.Item(30, 9).Style = "Normal"
$ws.Cells.Item(30, 10).NumberFormat = "@"
$ws.Cells.Item(30, 10).Value = "8.74"
$ws.Cells.Item(30, 10).Style = "Normal"
$ws.Cells.Item(30, 11).NumberFormat = "@"
$ws.Cells.Item(30, 11).Value = "8-32"
$ws.Cells.Item(30, 11).Style = "Normal"
$ws.Cells.Item(30, 12).NumberFormat = "@"
$ws.Cells.Item(30, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(30, 12).Style = "Normal"
# row 31
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "5/16`""
$ws.Cells.Item(31, 1).Style = "Normal"
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "Fully Threaded"
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "0.373`""
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.133`""
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "T25"
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(31, 6).NumberFormat = "@"
$ws.Cells.Item(31, 6).Value = "Not Rated"
$ws.Cells.Item(31, 6).Style = "Normal"
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(31, 8).NumberFormat = "@"
$ws.Cells.Item(31, 8).Value = "50"
$ws.Cells.Item(31, 8).Style = "Normal"
$ws.Cells.Item(31, 9).NumberFormat = "@"
$ws.Cells.Item(31, 9).Value = "97690A173"
$ws.Cells.Item(31, 9).Style = "Normal"
$ws.Cells.Item(31, 10).NumberFormat = "@"
$ws.Cells.Item(31, 10).Value = "10.71"
$ws.Cells.Item(31, 10).Style = "Normal"
$ws.Cells.Item(31, 11).NumberFormat = "@"
$ws.Cells.Item(31, 11).Value = "10-32"
$ws.Cells.Item(31, 11).Style = "Normal"
$ws.Cells.Item(31, 12).NumberFormat = "@"
$ws.Cells.Item(31, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(31, 12).Style = "Normal"
# row 32
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "3/8`""
$ws.Cells.Item(32, 1).Style = "Normal"
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = "Fully Threaded"
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "0.373`""
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.133`""
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "T25"
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(32, 6).NumberFormat = "@"
$ws.Cells.Item(32, 6).Value = "Not Rated"
$ws.Cells.Item(32, 6).Style = "Normal"
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(32, 8).NumberFormat = "@"
$ws.Cells.Item(32, 8).Value = "25"
$ws.Cells.Item(32, 8).Style = "Normal"
$ws.Cells.Item(32, 9).NumberFormat = "@"
$ws.Cells.Item(32, 9).Value = "97690A174"
$ws.Cells.Item(32, 9).Style = "Normal"
$ws.Cells.Item(32, 10).NumberFormat = "@"
$ws.Cells.Item(32, 10).Value = "6.60"
$ws.Cells.Item(32, 10).Style = "Normal"
$ws.Cells.Item(32, 11).NumberFormat = "@"
$ws.Cells.Item(32, 11).Value = "10-32"
$ws.Cells.Item(32, 11).Style = "Normal"
$ws.Cells.Item(32, 12).NumberFormat = "@"
$ws.Cells.Item(32, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(32, 12).Style = "Normal"
# row 33
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "1/2`""
$ws.Cells.Item(33, 1).Style = "Normal"
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = "Fully Threaded"
$ws.Cells.Item(33, 2).Style = "Normal"
$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = "0.373`""
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.133`""
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "T25"
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(33, 6).NumberFormat = "@"
$ws.Cells.Item(33, 6).Value = "Not Rated"
$ws.Cells.Item(33, 6).Style = "Normal"
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(33, 8).NumberFormat = "@"
$ws.Cells.Item(33, 8).Value = "25"
$ws.Cells.Item(33, 8).Style = "Normal"
$ws.Cells.Item(33, 9).NumberFormat = "@"
$ws.Cells.Item(33, 9).Value = "97690A175"
$ws.Cells.Item(33, 9).Style = "Normal"
$ws.Cells.Item(33, 10).NumberFormat = "@"
$ws.Cells.Item(33, 10).Value = "6.93"
$ws.Cells.Item(33, 10).Style = "Normal"
$ws.Cells.Item(33, 11).NumberFormat = "@"
$ws.Cells.Item(33, 11).Value = "10-32"
$ws.Cells.Item(33, 11).Style = "Normal"
$ws.Cells.Item(33, 12).NumberFormat = "@"
$ws.Cells.Item(33, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(33, 12).Style = "Normal"
# row 34
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "3/4`""
$ws.Cells.Item(34, 1).Style = "Normal"
$ws.Cells.Item(34, 2).NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = "Fully Threaded"
$ws.Cells.Item(34, 2).Style = "Normal"
$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = "0.373`""
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.133`""
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "T25"
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(34, 6).NumberFormat = "@"
$ws.Cells.Item(34, 6).Value = "Not Rated"
$ws.Cells.Item(34, 6).Style = "Normal"
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(34, 7).Style = "Normal"
$ws.Cells.Item(34, 8).NumberFormat = "@"
$ws.Cells.Item(34, 8).Value = "25"
$ws.Cells.Item(34, 8).Style = "Normal"
$ws.Cells.Item(34, 9).NumberFormat = "@"
$ws.Cells.Item(34, 9).Value = "97690A176"
$ws.Cells.Item(34, 9).Style = "Normal"
$ws.Cells.Item(34, 10).NumberFormat = "@"
$ws.Cells.Item(34, 10).Value = "8.22"
$ws.Cells.Item(34, 10).Style = "Normal"
$ws.Cells.Item(34, 11).NumberFormat = "@"
$ws.Cells.Item(34, 11).Value = "10-32"
$ws.Cells.Item(34, 11).Style = "Normal"
$ws.Cells.Item(34, 12).NumberFormat = "@"
$ws.Cells.Item(34, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(34, 12).Style = "Normal"
# row 35
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "1`""
$ws.Cells.Item(35, 1).Style = "Normal"
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = "Fully Threaded"
$ws.Cells.Item(35, 2).Style = "Normal"
$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = "0.373`""
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.133`""
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "T25"
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(35, 6).NumberFormat = "@"
$ws.Cells.Item(35, 6).Value = "Not Rated"
$ws.Cells.Item(35, 6).Style = "Normal"
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(35, 8).NumberFormat = "@"
$ws.Cells.Item(35, 8).Value = "10"
$ws.Cells.Item(35, 8).Style = "Normal"
$ws.Cells.Item(35, 9).NumberFormat = "@"
$ws.Cells.Item(35, 9).Value = "97690A177"
$ws.Cells.Item(35, 9).Style = "Normal"
$ws.Cells.Item(35, 10).NumberFormat = "@"
$ws.Cells.Item(35, 10).Value = "5.19"
$ws.Cells.Item(35, 10).Style = "Normal"
$ws.Cells.Item(35, 11).NumberFormat = "@"
$ws.Cells.Item(35, 11).Value = "10-32"
$ws.Cells.Item(35, 11).Style = "Normal"
$ws.Cells.Item(35, 12).NumberFormat = "@"
$ws.Cells.Item(35, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(35, 12).Style = "Normal"
# row 36
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "1/2`""
$ws.Cells.Item(36, 1).Style = "Normal"
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "Fully Threaded"
$ws.Cells.Item(36, 2).Style = "Normal"
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = "0.492`""
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.175`""
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "T30"
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(36, 6).NumberFormat = "@"
$ws.Cells.Item(36, 6).Value = "Not Rated"
$ws.Cells.Item(36, 6).Style = "Normal"
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value = "10"
$ws.Cells.Item(36, 8).Style = "Normal"
$ws.Cells.Item(36, 9).NumberFormat = "@"
$ws.Cells.Item(36, 9).Value = "97690A178"
$ws.Cells.Item(36, 9).Style = "Normal"
$ws.Cells.Item(36, 10).NumberFormat = "@"
$ws.Cells.Item(36, 10).Value = "5.60"
$ws.Cells.Item(36, 10).Style = "Normal"
$ws.Cells.Item(36, 11).NumberFormat = "@"
$ws.Cells.Item(36, 11).Value = "1/4`"-20"
$ws.Cells.Item(36, 11).Style = "Normal"
$ws.Cells.Item(36, 12).NumberFormat = "@"
$ws.Cells.Item(36, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(36, 12).Style = "Normal"
# row 37
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "3/4`""
$ws.Cells.Item(37, 1).Style = "Normal"
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = "Fully Threaded"
$ws.Cells.Item(37, 2).Style = "Normal"
$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "0.492`""
$ws.Cells.Item(37, 3).Style = "Normal"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.175`""
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "T30"
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(37, 6).NumberFormat = "@"
$ws.Cells.Item(37, 6).Value = "Not Rated"
$ws.Cells.Item(37, 6).Style = "Normal"
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value = "10"
$ws.Cells.Item(37, 8).Style = "Normal"
$ws.Cells.Item(37, 9).NumberFormat = "@"
$ws.Cells.Item(37, 9).Value = "97690A179"
$ws.Cells.Item(37, 9).Style = "Normal"
$ws.Cells.Item(37, 10).NumberFormat = "@"
$ws.Cells.Item(37, 10).Value = "6.49"
$ws.Cells.Item(37, 10).Style = "Normal"
$ws.Cells.Item(37, 11).NumberFormat = "@"
$ws.Cells.Item(37, 11).Value = "1/4`"-20"
$ws.Cells.Item(37, 11).Style = "Normal"
$ws.Cells.Item(37, 12).NumberFormat = "@"
$ws.Cells.Item(37, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(37, 12).Style = "Normal"
# row 38
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "1`""
$ws.Cells.Item(38, 1).Style = "Normal"
$ws.Cells.Item(38, 2).NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = "Fully Threaded"
$ws.Cells.Item(38, 2).Style = "Normal"
$ws.Cells.Item(38, 3).NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = "0.492`""
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.175`""
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "T30"
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(38, 6).NumberFormat = "@"
$ws.Cells.Item(38, 6).Value = "Not Rated"
$ws.Cells.Item(38, 6).Style = "Normal"
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(38, 8).NumberFormat = "@"
$ws.Cells.Item(38, 8).Value = "10"
$ws.Cells.Item(38, 8).Style = "Normal"
$ws.Cells.Item(38, 9).NumberFormat = "@"
$ws.Cells.Item(38, 9).Value = "97690A181"
$ws.Cells.Item(38, 9).Style = "Normal"
$ws.Cells.Item(38, 10).NumberFormat = "@"
$ws.Cells.Item(38, 10).Value = "7.76"
$ws.Cells.Item(38, 10).Style = "Normal"
$ws.Cells.Item(38, 11).NumberFormat = "@"
$ws.Cells.Item(38, 11).Value = "1/4`"-20"
$ws.Cells.Item(38, 11).Style = "Normal"
$ws.Cells.Item(38, 12).NumberFormat = "@"
$ws.Cells.Item(38, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(38, 12).Style = "Normal"
# row 39
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "1 1/4`""
$ws.Cells.Item(39, 1).Style = "Normal"
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = "Fully Threaded"
$ws.Cells.Item(39, 2).Style = "Normal"
$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = "0.492`""
$ws.Cells.Item(39, 3).Style = "Normal"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.175`""
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "T30"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 6).NumberFormat = "@"
$ws.Cells.Item(39, 6).Value = "Not Rated"
$ws.Cells.Item(39, 6).Style = "Normal"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "ASME B18.6.3"
$ws.Cells.Item(39, 7).Style = "Normal"
$ws.Cells.Item(39, 8).NumberFormat = "@"
$ws.Cells.Item(39, 8).Value = "10"
$ws.Cells.Item(39, 8).Style = "Normal"
$ws.Cells.Item(39, 9).NumberFormat = "@"
$ws.Cells.Item(39, 9).Value = "97690A182"
$ws.Cells.Item(39, 9).Style = "Normal"
$ws.Cells.Item(39, 10).NumberFormat = "@"
$ws.Cells.Item(39, 10).Value = "9.58"
$ws.Cells.Item(39, 10).Style = "Normal"
$ws.Cells.Item(39, 11).NumberFormat = "@"
$ws.Cells.Item(39, 11).Value = "1/4`"-20"
$ws.Cells.Item(39, 11).Style = "Normal"
$ws.Cells.Item(39, 12).NumberFormat = "@"
$ws.Cells.Item(39, 12).Value = "Zinc-Plated Steel"
$ws.Cells.Item(39, 12).Style = "Normal"
